# msz - floProfil first tc is running
# Fill in the first test-case row (row 3) of the profile-dialog test sheet
# with a first name / last name, then leave the selection on A3.
#
# NOTE on ordering: the shared-strings table is appended-to in the order
# new unique strings are first encountered, so "Schmotz" (I3) is written
# before "Matthias" (H3) to reproduce the author's sharedStrings.xml order
# (index 24 = Schmotz, index 25 = Matthias) and the resulting <c> value
# indices (H3 -> 25, I3 -> 24).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I3").Value = "Schmotz"
$ws.Range("H3").Value = "Matthias"

[void]$ws.Range("A3").Select()
